$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# Row 13 (Tags row): swap "DNASeq" out of column C, move "Assembly" into C,
# and add the new "DNA Sequencing" tag in column D.
$ws.Range("C13").Value = "Assembly"
$ws.Range("D13").Value = "DNA Sequencing"

# Row 14 (Tags Term Accession Number row): shift the existing NCIT_C52474
# accession into column C (aligning with "Assembly"), and add the new
# accession number for "DNA Sequencing" in column D.
$ws.Range("C14").Value = "http://purl.obolibrary.org/obo/NCIT_C52474"
$ws.Range("D14").Value = "http://purl.obolibrary.org/obo/NCIT_C153598"

# Row 15 (Tags Term Source REF row): add the ontology source ref for the
# newly added "DNA Sequencing" tag in column D.
$ws.Range("D15").Value = "NCIT"
